# Auto-generated edit script: Module: Meeting Notes
$wb = $excel.ActiveWorkbook

# --- Firm sheet: add record row 51 ---
$wsFirm = $wb.Worksheets.Item("Firm")
$wsFirm.Range("A51").Value = "AMNNR_Record003"
$wsFirm.Range("B51").Value = "Acc 3"
$wsFirm.Range("C51").Value = "Company"

# --- Contact sheet: add Meeting Notes Notification Reminder section ---
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Range("A36:I36").Value = "Meeting Notes Notification Reminder"
$wsContact.Range("A36:I36").Merge()
$wsAcuityRef = $wb.Worksheets.Item("Acuity")
$wsAcuityRef.Range("A70:I70").Copy()
$wsContact.Range("A36:I36").PasteSpecial(-4122)
$wsContact.Range("A36").Value = "Meeting Notes Notification Reminder"
$wsContact.Range("A37").Value = "AMNNR_Con1"
$wsContact.Range("D37").Value = "con 5"
$wsContact.Range("A38").Value = "AMNNR_Con2"
$wsContact.Range("D38").Value = "Jhon"
$wsContact.Range("A39").Value = "AMNNR_Con3"
$wsContact.Range("D39").Value = "Max"
$wsContact.Range("A40").Value = "AMNNR_Con4"
$wsContact.Range("D40").Value = "Martha"

# --- Activity Timeline sheet: add Meeting Notes Notification Reminder activity rows ---
$wsActivity = $wb.Worksheets.Item("Activity Timeline")
$wsActivity.Range("A125").Value = "AMNNR_Activity006"
$wsActivity.Range("C125").Value = "Introduction"
$wsActivity.Range("E125").Value = "con 5<break>con 6<break>Sumo Logic<break>Vertica<break>Demo Deal<break>Mutual Fund"
$wsActivity.Range("M125").Value = "Completed"
$wsActivity.Range("N125").Value = "Normal"
$wsActivity.Range("A126").Value = "AMNNR_Activity007"
$wsActivity.Range("D126").Value = "areca  moss fundraising should be tagged"
$wsActivity.Range("G126").Value = "All Records Select"
$wsActivity.Range("A127").Value = "AMNNR_Activity008"
$wsActivity.Range("C127").Value = "Sales Meeting"
$wsActivity.Range("E127").Value = "Jhon<break>con 10<break>Sumo Logic"
$wsActivity.Range("M127").Value = "Completed"
$wsActivity.Range("N127").Value = "Normal"
$wsActivity.Range("A128").Value = "AMNNR_Activity009"
$wsActivity.Range("D128").Value = "Golden Ret"
$wsActivity.Range("G128").Value = "Golden Ret"
$wsActivity.Range("A129").Value = "AMNNR_Activity010"
$wsActivity.Range("C129").Value = "Marketing Strategy"
$wsActivity.Range("D129").Value = "We as an organization need to have certain strategy towards our marketing approch with Vertica and sumo logic Firm"
$wsActivity.Range("E129").Value = "Max<break>Jhon<break>con 11"
$wsActivity.Range("G129").Value = "Vertica<break>Sumo Logic"
$wsActivity.Range("M129").Value = "Not Started"
$wsActivity.Range("N129").Value = "Normal"
$wsActivity.Range("A130").Value = "AMNNR_Activity011"
$wsActivity.Range("D130").Value = "keep in loop  con 4, con 5"
$wsActivity.Range("E130").Value = "<break>Vertica<break>Sumo Logic"
$wsActivity.Range("G130").Value = "con 4<break>con 5"
$wsActivity.Range("A131").Value = "AMNNR_Activity012"
$wsActivity.Range("E131").Value = "<break>con 4<break>con 5"
$wsActivity.Range("A132").Value = "AMNNR_Activity013"
$wsActivity.Range("C132").Value = "Send Quote"
$wsActivity.Range("D132").Value = "unicorn"
$wsActivity.Range("E132").Value = "Maxtra<break>Martha<break>Jhon<break>con 11<break>Sumo Logic<break>Vertica"
$wsActivity.Range("M132").Value = "In Progress"
$wsActivity.Range("N132").Value = "Normal"
$wsActivity.Range("A133").Value = "AMNNR_Activity014"
$wsActivity.Range("D133").Value = "Palm areca"
$wsActivity.Range("E133").Value = "<break>Palm<break>areca"
$wsActivity.Range("G133").Value = "All Records Select"
$wsActivity.Range("A134").Value = "AMNNR_Activity015"
$wsActivity.Range("C134").Value = "Send Notice"
$wsActivity.Range("E134").Value = "Acc 3<break>Martha<break>Echo<break>Alexa<break>Green Pothos<break>areca"
$wsActivity.Range("M134").Value = "In Progress"
$wsActivity.Range("N134").Value = "Normal"
$wsActivity.Range("A135").Value = "AMNNR_Activity016"
$wsActivity.Range("C135").Value = "Send Notice updated"
$wsActivity.Range("A136").Value = "AMNNR_Activity017"
$wsActivity.Range("C136").Value = "SSend Notice"
$wsActivity.Range("D136").Value = "echo alexa Green pothos areca"
$wsActivity.Range("E136").Value = "Acc 3<break>Martha"
$wsActivity.Range("M136").Value = "In Progress"
$wsActivity.Range("N136").Value = "Normal"
$wsActivity.Range("A137").Value = "AMNNR_Activity018"
$wsActivity.Range("C137").Value = "SSend Notice Follow up 1"
$wsActivity.Range("A138").Value = "AMNNR_Activity019"
$wsActivity.Range("C138").Value = "SSend Notice Follow up 2"
$wsActivity.Range("A139").Value = "AMNNR_Activity020"
$wsActivity.Range("C139").Value = "SSend Notice Follow up 3"
$wsActivity.Range("A140").Value = "AMNNR_Activity021"
$wsActivity.Range("D140").Value = "Follow up task As Send Notice Updated for Con 6,Con 7, Con 8, Acc 4"
$wsActivity.Range("E140").Value = "<break>con 6<break>con 7<break>con 8<break>Acc 4"
$wsActivity.Range("G140").Value = "con 6<break>con 7<break>con 8<break>Acc 4"
$wsActivity.Range("A141").Value = "AMNNR_Activity022"
$wsActivity.Range("E141").Value = "Mutual Fund<break>FC Fundraising<break>Acc 1"
$wsActivity.Range("A142").Value = "AMNNR_Activity023"
$wsActivity.Range("C142").Value = "Task for the day"
$wsActivity.Range("D142").Value = "Follow up with Contacts Con 4, Con 5 about demo deal"
$wsActivity.Range("E142").Value = "Con 1<break>con 2<break>Acc 3<break>Maxtra"
$wsActivity.Range("G142").Value = "con 4<break>con 5<break>Demo Deal"
$wsActivity.Range("M142").Value = "In Progress"
$wsActivity.Range("N142").Value = "Normal"
$wsActivity.Range("A143").Value = "AMNNR_Activity024"
$wsActivity.Range("E143").Value = "<break>con 4<break>con 5<break>Demo Deal"
$wsActivity.Range("A144").Value = "AMNNR_Activity025"
$wsActivity.Range("E144").Value = "Maxtra<break>Demo Deal"
$wsActivity.Range("A145").Value = "AMNNR_Activity026"
$wsActivity.Range("E145").Value = "Con 1<break>con 2<break>Acc 3<break>con 4<break>con 5"
$wsActivity.Range("A146").Value = "AMNNR_Activity027"
$wsActivity.Range("C146").Value = "Task Test"
$wsActivity.Range("D146").Value = "Follow up with Contacts Con 4, Con 5 about demo deal"
$wsActivity.Range("E146").Value = "Con 1<break>con 2<break>Acc 3<break>Maxtra"
$wsActivity.Range("G146").Value = "con 4<break>con 5<break>Demo Deal"
$wsActivity.Range("M146").Value = "In Progress"
$wsActivity.Range("N146").Value = "Normal"
$wsActivity.Range("A147").Value = "AMNNR_Activity028"
$wsActivity.Range("E147").Value = "<break>con 4<break>con 5<break>Demo Deal"
$wsActivity.Range("A148").Value = "AMNNR_Activity029"
$wsActivity.Range("E148").Value = "Contact Invalid<break>Account Invalid"

# Apply wrap-text style to Notes (D) column cells, matching existing style used in column D
$wsActivity.Range("D124").Copy()
$wsActivity.Range("D126").PasteSpecial(-4122)
$wsActivity.Range("D128").PasteSpecial(-4122)
$wsActivity.Range("D129").PasteSpecial(-4122)
$wsActivity.Range("D130").PasteSpecial(-4122)
$wsActivity.Range("D132").PasteSpecial(-4122)
$wsActivity.Range("D133").PasteSpecial(-4122)
$wsActivity.Range("D136").PasteSpecial(-4122)
$wsActivity.Range("D140").PasteSpecial(-4122)
$wsActivity.Range("D142").PasteSpecial(-4122)
$wsActivity.Range("D146").PasteSpecial(-4122)
$wsActivity.Range("D126").Value = "areca  moss fundraising should be tagged"
$wsActivity.Range("D128").Value = "Golden Ret"
$wsActivity.Range("D129").Value = "We as an organization need to have certain strategy towards our marketing approch with Vertica and sumo logic Firm"
$wsActivity.Range("D130").Value = "keep in loop  con 4, con 5"
$wsActivity.Range("D132").Value = "unicorn"
$wsActivity.Range("D133").Value = "Palm areca"
$wsActivity.Range("D136").Value = "echo alexa Green pothos areca"
$wsActivity.Range("D140").Value = "Follow up task As Send Notice Updated for Con 6,Con 7, Con 8, Acc 4"
$wsActivity.Range("D142").Value = "Follow up with Contacts Con 4, Con 5 about demo deal"
$wsActivity.Range("D146").Value = "Follow up with Contacts Con 4, Con 5 about demo deal"

# --- Acuity sheet: add Meeting Notes Notification Reminder connection rows ---
$wsAcuity = $wb.Worksheets.Item("Acuity")
$wsAcuity.Range("A74").Value = "AMNNR_Acuity004"
$wsAcuity.Range("AA74").Value = "con 5<break>con 6<break>+5"
$wsAcuity.Range("A75").Value = "AMNNR_Acuity005"
$wsAcuity.Range("AA75").Value = "con 5<break>con 6<break>+10"
$wsAcuity.Range("A76").Value = "AMNNR_Acuity006"
$wsAcuity.Range("AA76").Value = "Jhon<break>con 10<break>+2"
$wsAcuity.Range("A77").Value = "AMNNR_Acuity007"
$wsAcuity.Range("AA77").Value = "Jhon<break>con 10<break>+3"
$wsAcuity.Range("A78").Value = "AMNNR_Acuity008"
$wsAcuity.Range("AA78").Value = "Max<break>Jhon<break>+4"
$wsAcuity.Range("A79").Value = "AMNNR_Acuity009"
$wsAcuity.Range("AA79").Value = "Max<break>Jhon<break>+6"
$wsAcuity.Range("A80").Value = "AMNNR_Acuity010"
$wsAcuity.Range("AA80").Value = "Martha<break>Jhon<break>+5"
$wsAcuity.Range("A81").Value = "AMNNR_Acuity011"
$wsAcuity.Range("AA81").Value = "Martha<break>Jhon<break>+7"
$wsAcuity.Range("A82").Value = "AMNNR_Acuity012"
$wsAcuity.Range("AA82").Value = "Martha<break>areca<break>+5"
$wsAcuity.Range("A83").Value = "AMNNR_Acuity013"
$wsAcuity.Range("AA83").Value = "<break>Martha<break>+1"
$wsAcuity.Range("A84").Value = "AMNNR_Acuity014"
$wsAcuity.Range("AA84").Value = "Martha<break>con 6<break>+5"
$wsAcuity.Range("A85").Value = "AMNNR_Acuity015"
$wsAcuity.Range("AA85").Value = "Martha<break>con 6<break>+8"
$wsAcuity.Range("A86").Value = "AMNNR_Acuity016"
$wsAcuity.Range("AA86").Value = "Con 1<break>con 2<break>+6"
$wsAcuity.Range("A87").Value = "AMNNR_Acuity017"
$wsAcuity.Range("AA87").Value = "Con 1<break>con 2<break>+4"
$wsAcuity.Range("A88").Value = "AMNNR_Acuity018"
$wsAcuity.Range("AA88").Value = "Con 1<break>con 2<break>+6"

# --- Update view/selection state to match final workbook state ---
$wsFirm.Activate()
$wsFirm.Range("B55").Select()

$wsContact.Activate()
$wsContact.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsContact.Range("A40").Select()

$wsActivity.Activate()
$wsActivity.Range("G146").Select()

$wsAcuity.Activate()
$wsAcuity.Range("B88").Select()

